$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a fresh row at position 47. This shifts the existing ---
# --- rows 47 ("Femacal de La Calera" / 2022-07-04 entry) and 48 (2021-05-28 ---
# --- entry) down to rows 48 and 49 respectively, preserving their values. ---
$ws.Rows.Item(47).Insert()

# --- Step 2: fill the newly-inserted row 47 with the new weekly price entry ---
# --- (2022-07-06), re-using the static descriptive columns from the row ---
# --- immediately below (now row 48, the former row 47). ---
$ws.Range("A47").Value = 3
$ws.Range("B47").Value = "Femacal de La Calera"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 44748
$ws.Range("E47").Value = 5
$ws.Range("F47").Value = 100112035
$ws.Range("G47").Value = "Bruselas (repollito)"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 73
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 16000
$ws.Range("M47").Value = 15521
$ws.Range("N47").Value = "$/malla 15 kilos"
$ws.Range("O47").Value = "Provincia de Quillota"
$ws.Range("P47").Value = 1035
$ws.Range("Q47").Value = 15
$ws.Range("R47").Value = "Hortaliza"

# --- Step 3: append a brand-new row 50 with the remaining new weekly entry ---
# --- (2022-07-05), following the last existing row (49). ---
$ws.Range("A50").Value = 3
$ws.Range("B50").Value = "Femacal de La Calera"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 44747
$ws.Range("D50").NumberFormat = $ws.Range("D49").NumberFormat
$ws.Range("E50").Value = 5
$ws.Range("F50").Value = 100112035
$ws.Range("G50").Value = "Bruselas (repollito)"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 40
$ws.Range("K50").Value = 16000
$ws.Range("L50").Value = 16000
$ws.Range("M50").Value = 16000
$ws.Range("N50").Value = "$/malla 15 kilos"
$ws.Range("O50").Value = "Provincia de Quillota"
$ws.Range("P50").Value = 1067
$ws.Range("Q50").Value = 15
$ws.Range("R50").Value = "Hortaliza"
